# Release version 3.41 (#170)
# Wrap the QR-code placeholder ("{%DemisIdQrImage}") and its "Meldungs-ID"
# caption - both living inside the VML text box anchored to the report
# header paragraph - with a DocxTemplater-style section:
#   {#DemisIdQrImage} ... {%DemisIdQrImage} ... Meldungs-ID ... {/DemisIdQrImage}
# Word's Find/Replace cannot reach into legacy VML <w:pict> text-box content,
# so the whole host paragraph is rebuilt via Range.InsertXML with the target
# WordprocessingML fragment (identical to the original except for the new
# runs plus the regenerated shape id / w14:anchorId that Word assigns when
# the drawing is re-serialized).

$d = $word.ActiveDocument

$targetParagraph = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Text -like "*Untersuchungsbefund*zu KL*") {
        $targetParagraph = $candidate
        break
    }
}

if ($targetParagraph -eq $null) {
    throw "Could not locate the 'Untersuchungsbefund' heading paragraph that hosts the QR code text box."
}

$newParagraphXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:o="urn:schemas-microsoft-com:office:office" w14:paraId="301E55E0" w14:textId="77777777" w:rsidR="005B2384" w:rsidRPr="003136E9" w:rsidRDefault="005B2384" w:rsidP="009B6F20"><w:pPr><w:widowControl w:val="0"/><w:tabs><w:tab w:val="left" w:pos="90"/><w:tab w:val="left" w:pos="3231"/></w:tabs><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:before="240" w:after="160"/><w:ind w:right="2381"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:bCs/><w:color w:val="000000"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:noProof/></w:rPr><w:pict w14:anchorId="67377B5C"><v:shapetype id="_x0000_t202" coordsize="21600,21600" o:spt="202" path="m,l,21600r21600,l21600,xe"><v:stroke joinstyle="miter"/><v:path gradientshapeok="t" o:connecttype="rect"/></v:shapetype><v:shape id="_x0000_s2052" type="#_x0000_t202" style="position:absolute;margin-left:249.4pt;margin-top:47pt;width:90.7pt;height:99.2pt;z-index:-1;visibility:visible;mso-wrap-distance-left:9pt;mso-wrap-distance-top:3.6pt;mso-wrap-distance-right:9pt;mso-wrap-distance-bottom:3.6pt;mso-position-horizontal-relative:text;mso-position-vertical-relative:text;mso-width-relative:margin;mso-height-relative:margin;v-text-anchor:top" stroked="f"><v:textbox><w:txbxContent><w:p w14:paraId="2B738C91" w14:textId="77777777" w:rsidR="005B2384" w:rsidRDefault="005B2384" w:rsidP="00EA2A76"><w:pPr><w:jc w:val="center"/><w:rPr><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="12"/><w:szCs w:val="12"/><w:lang w:val="en-US"/></w:rPr><w:t>{#</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="12"/><w:szCs w:val="12"/><w:lang w:val="en-US"/></w:rPr><w:t>DemisIdQrImage</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="12"/><w:szCs w:val="12"/><w:lang w:val="en-US"/></w:rPr><w:t>}</w:t></w:r><w:r w:rsidRPr="006F7F3C"><w:rPr><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t>{%DemisIdQrImage}</w:t></w:r></w:p><w:p w14:paraId="0AE45C70" w14:textId="77777777" w:rsidR="005B2384" w:rsidRPr="000E6D63" w:rsidRDefault="005B2384" w:rsidP="00EA2A76"><w:pPr><w:jc w:val="center"/><w:rPr><w:sz w:val="12"/><w:szCs w:val="12"/></w:rPr></w:pPr><w:r w:rsidRPr="000E6D63"><w:rPr><w:sz w:val="12"/><w:szCs w:val="12"/></w:rPr><w:t>Meldungs-ID</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="12"/><w:szCs w:val="12"/><w:lang w:val="en-US"/></w:rPr><w:t>{</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="12"/><w:szCs w:val="12"/><w:lang w:val="en-US"/></w:rPr><w:t>/</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="12"/><w:szCs w:val="12"/><w:lang w:val="en-US"/></w:rPr><w:t>DemisIdQrImage</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="12"/><w:szCs w:val="12"/><w:lang w:val="en-US"/></w:rPr><w:t>}</w:t></w:r></w:p></w:txbxContent></v:textbox></v:shape></w:pict></w:r><w:r w:rsidRPr="00282FF8"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:color w:val="000000"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve">Untersuchungsbefund </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:color w:val="000000"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>zu KL</w:t></w:r><w:r w:rsidRPr="00CA336F"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:color w:val="000000"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>{LaboratoryNumber}</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:color w:val="000000"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidRPr="00282FF8"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:color w:val="000000"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>-</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:color w:val="000000"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidRPr="00282FF8"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:i/><w:iCs/><w:color w:val="000000"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>Endbefund</w:t></w:r></w:p>
'@

[void]$targetParagraph.Range.InsertXML($newParagraphXml)

Write-Output "Updated QR code text box placeholders."
